$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.679.78'
$ws.Range('E2').Value = '  -2.54%  '
$ws.Range('D3').Value = '1.556.74'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''205.84'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('E6').Value = '  -2.08%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '''21.93'
$ws.Range('D8').ClearFormats()
$ws.Range('E9').Value = '  -0.56%  '
$ws.Range('E10').Value = '  -1.53%  '
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('D12').Value = '1.777.71'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').Value = '1.561.61'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').Value = '''3.74'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('D15').Value = '''0.512'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').Value = '''61.65'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.76%  '
$ws.Range('D17').Value = '26.730.74'
$ws.Range('E17').Value = '  -2.39%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D18').Value = '''213.47'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''7.33'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.95%  '
$ws.Range('E20').Value = '  -2.02%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('D23').Value = '''9.35'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.77%  '
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('D26').Value = '''6.77'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('D27').Value = '''14.81'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.13%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('D30').Value = '''0.0463'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.51%  '
$ws.Range('E31').Value = '  -3.70%  '
$ws.Range('D32').Value = '''3.15'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('D33').Value = '1.385.48'
$ws.Range('E33').Value = '  +1.57%  '
$ws.Range('E34').Value = '  -1.57%  '
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('E36').Value = '  -0.94%  '
$ws.Range('E37').Value = '  -4.16%  '
$ws.Range('D38').Value = '''0.0163'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.46%  '
$ws.Range('D39').Value = '''0.517'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.82%  '
$ws.Range('D40').Value = '''0.811'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('D42').Value = '''0.994'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.21%  '
$ws.Range('E43').Value = '  +1.88%  '
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('E45').Value = '  -2.07%  '
$ws.Range('D46').Value = '''63.06'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.48%  '
$ws.Range('D47').Value = '1.691.72'
$ws.Range('D48').Value = '''85.35'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').Value = '0.0₇0973'
$ws.Range('E49').Value = '  -1.77%  '
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').Value = '''0.0947'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.74%  '
